# Refresh the crypto price/volume table (rows 2-51) with the latest scraped
# values.  Most cells just get new text; a handful of Price cells in column D
# now look like plain decimal numbers (single '.') which Excel would happily
# reinterpret as numeric values, so those are forced to Text format first to
# keep them stored as strings, matching the original inline-string content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# addr -> new value
$updates = [ordered]@{
    "D2"  = "68.309.08"
    "E2"  = "  -1.76%  "
    "D3"  = "2.445.51"
    "E3"  = "  -1.38%  "
    "E4"  = "  +0.00%  "
    "D5"  = "555.88"
    "E5"  = "  -2.09%  "
    "D6"  = "162.40"
    "E6"  = "  -1.67%  "
    "E7"  = "  -0.01%  "
    "D8"  = "0.499"
    "E8"  = "  -2.24%  "
    "D9"  = "2.445.01"
    "E9"  = "  -1.33%  "
    "E10" = "  -5.95%  "
    "E11" = "  -1.90%  "
    "E12" = "  -5.36%  "
    "D13" = "4.74"
    "E13" = "  -3.29%  "
    "D14" = "2.895.66"
    "E14" = "  -1.37%  "
    "D15" = "68.224.22"
    "E15" = "  -1.89%  "
    "D16" = "0.0000166"
    "E16" = "  -4.48%  "
    "D17" = "23.03"
    "E17" = "  -4.97%  "
    "D18" = "2.455.85"
    "E18" = "  -1.20%  "
    "D19" = "10.79"
    "E19" = "  -2.98%  "
    "D20" = "338.21"
    "E20" = "  -1.89%  "
    "D21" = "7.07"
    "E21" = "  -3.81%  "
    "D22" = "3.72"
    "E22" = "  -3.37%  "
    "D23" = "1.00"
    "E23" = "  +0.01%  "
    "E24" = "  -4.67%  "
    "D25" = "67.10"
    "E25" = "  -4.49%  "
    # Rows 26/27 swap places (NEARProtocol <-> WrappedeETH) and get new figures
    "B26" = "WrappedeETH"
    "C26" = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
    "D26" = "2.573.39"
    "E26" = "  -1.60%  "
    "B27" = "NEARProtocol"
    "C27" = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
    "D27" = "3.63"
    "E27" = "  -6.47%  "
    "D28" = "0.998"
    "E28" = "  +0.07%  "
    "D29" = "8.01"
    "E29" = "  -7.13%  "
    "D30" = "0.0₃0817"
    "E30" = "  -6.58%  "
    "D31" = "7.08"
    "E31" = "  -8.87%  "
    "D32" = "1.00"
    "E32" = "  +0.04%  "
    "D33" = "423.25"
    "E33" = "  -4.43%  "
    "E34" = "  -4.16%  "
    "E35" = "  -4.48%  "
    "D36" = "157.13"
    "E36" = "  +0.90%  "
    "E37" = "  -0.20%  "
    "E38" = "  +0.05%  "
    "E39" = "  -4.27%  "
    "D40" = "17.65"
    "E40" = "  -2.60%  "
    "D41" = "0.300"
    "E41" = "  -4.29%  "
    "D42" = "4.36"
    "E42" = "  -4.87%  "
    "E43" = "  -6.26%  "
    "E44" = "  +0.16%  "
    "D45" = "133.28"
    "E45" = "  -4.32%  "
    "D46" = "2.03"
    "E46" = "  -6.43%  "
    "D47" = "3.30"
    "E47" = "  -3.77%  "
    "D48" = "0.0712"
    "E48" = "  -2.24%  "
    "E49" = "  -7.11%  "
    "E50" = "  -2.66%  "
    "D51" = "0.0902"
    "E51" = "  -1.79%  "
}

foreach ($addr in $updates.Keys) {
    $newValue = $updates[$addr]
    $cell = $ws.Range($addr)

    # Column D holds price strings like "68.309.08" (thousand separators) that
    # are not valid numbers, but some new values (e.g. "555.88", "1.00") are
    # syntactically valid numbers. Force those specific cells to Text format
    # *before* writing so Excel keeps storing them as plain strings, just
    # like the rest of the (inline-string) price column.
    if ($addr.StartsWith("D") -and ($newValue -match '^[+-]?\d+(\.\d+)?$')) {
        $cell.NumberFormat = "@"
    }

    $cell.Value = $newValue
}
